# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (a copy of the "2021-Q4" fund-holding
# layout, populated with the new quarter's figures) positioned between the
# existing "2021-Q4" and "总计" sheets, and prepends a matching row to the
# "总计" summary sheet.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1. Create the new "2022-Q1" sheet right after "2021-Q4" -------------
# (NOTE: worksheet references returned by Item(...) are index-bound here,
#  not stable object identities -- any sheet handle fetched *before* this
#  insertion, e.g. "总计", would silently resolve to the wrong tab
#  afterwards because everything after the insertion point shifts by one.
#  So "总计" is (re)fetched fresh, below, only once the sheet collection
#  has reached its final shape.)
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Clone the whole header+data block (values *and* formatting) from
# "2021-Q4" so fonts/borders/alignment match exactly, then overwrite just
# the cells whose values actually changed for the new quarter.
# (Row 1 has no A1 cell in the source sheet, so copy B:H for the header
#  separately from the A2 index cell, rather than a blanket A1:H2 which
#  would materialise a spurious empty A1.)
$q4.Range("B1:H1").Copy($q1.Range("B1:H1"))
$q4.Range("A2:H2").Copy($q1.Range("A2:H2"))

# Fund code/name are unchanged between quarters.
# Fund size / stock position / position share / market value are new --
# they are numeric-looking text in the source data (e.g. "6.05"), so each
# cell is briefly forced to Text format before assignment (otherwise
# Excel's normal type-inference silently turns them into numbers), then
# restored to the plain/default "Normal" style so only the value -- not
# the formatting -- differs from a freshly-cloned cell.
$q1.Cells.Item(2, 4).NumberFormat = "@"
$q1.Cells.Item(2, 4).Value = "6.05"
$q1.Cells.Item(2, 4).Style = "Normal"
$q1.Cells.Item(2, 5).NumberFormat = "@"
$q1.Cells.Item(2, 5).Value = "99.49"
$q1.Cells.Item(2, 5).Style = "Normal"
$q1.Cells.Item(2, 6).NumberFormat = "@"
$q1.Cells.Item(2, 6).Value = "8.51"
$q1.Cells.Item(2, 6).Style = "Normal"
$q1.Cells.Item(2, 7).NumberFormat = "@"
$q1.Cells.Item(2, 7).Value = "0.5149"
$q1.Cells.Item(2, 7).Style = "Normal"

# Position rank is a genuine number, and happens to be unchanged (2), so
# nothing further to do for column H.

# --- 2. Prepend a "2022-Q1" row in "总计" ---------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing data row down from row 2 to row 3 (value + format),
# then blank out row 2's contents (keeping A2's index-column formatting)
# ready for the new quarter's figures.
$total.Range("A2:D2").Copy($total.Range("A3:D3"))
$total.Range("A2:D2").ClearContents()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(3, 1).Value = 1

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.51
